$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Figure out the last used row (data rows below the header).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Insert two new columns before column B - this shifts the two existing
# week columns (B, C) two slots to the right, to D and E, carrying their
# values/styles/formatting along with them (matches the highlight style
# on the old C17 ending up on E17, and the old col min="3" width 8.0
# definition shifting out to col min="5").
$ws.Range("B1").EntireColumn.Insert()
$ws.Range("B1").EntireColumn.Insert()

# Give the two freshly inserted columns the same width (8 characters) as
# the shifted-out original data columns.
$ws.Columns("B").ColumnWidth = 7.14
$ws.Columns("C").ColumnWidth = 7.14

# New header labels for the two newly-inserted week columns.
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# Fill the new B/C columns for every data row with the same
# placeholder value that was carried into column D ("UN").
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($r, 4).Text
    $ws.Cells.Item($r, 3).Value = $ws.Cells.Item($r, 4).Text
}
